# working_hours.xlsx -- "bound all algorithm settings which were left"
#
# The sheet tracks working hours. Row 47 used to be an unused blank filler
# row directly above the summary block (sum [min] / sum [h] / sum [working
# weeks] in rows 48-50). This edit records one more working-hours entry in
# that row (2014-03-04, 08:40 -> 12:00) and pushes the blank filler row (and
# the summary block beneath it) down by one row, so the filler row is now 48
# and the summary rows are 49-51. The running SUM()/ratios are recalculated
# automatically by the formula engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 47. This pushes the old blank filler row (47) and the
# summary rows below it (48, 49, 50) down to 48, 49, 50, 51 respectively,
# and Excel automatically updates F49's SUM(F2:F47) -> SUM(F2:F48) reference
# because the inserted row lies inside that range.
$ws.Rows(47).Insert()

# Populate the freshly-inserted row 47 with the new working-hours entry.
$ws.Range("A47").Value = 2014
$ws.Range("B47").Value = 3
$ws.Range("C47").Value = 4
$ws.Range("D47").Value = 0.3611111111111111
$ws.Range("E47").Value = 0.5
$ws.Range("F47").Formula = "=(E47-D47)*24*60"
$ws.Range("G47").Formula = "=F47/60"

# Match the author's final cursor position/scroll in the diff.
$ws.Range("F47").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
